# Week 25 newsletter update — Rockies injuries sheet.
# Final layout (row 1 = header, unchanged):
#   Row 2: Nolan Arenado   / arenano01 / September 24 2017 / Hand     / bruised hand note
#   Row 3: Jairo Diaz      / diazja01  / September 12 2017 / Elbow    / elbow inflammation note (updated)
#   Row 4: Carlos Gonzalez / gonzaca01 / September 24 2017 / Shoulder / shoulder note
#   Row 5/6: blank placeholder rows, wrap-text style carried in column E only

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Nolan Arenado (new) ---
$ws.Range("A2").Value = "Nolan Arenado"
$ws.Range("B2").Value = "arenano01"
$ws.Range("C2").Value = "September 24 2017"
$ws.Range("D2").Value = "Hand"
$ws.Range("E2").Value = "Arenado is dealing with a bruised right hand near his thumb and his status for Monday's contest against the Marlins is undetermined."

# --- Row 3: Jairo Diaz (existing row, content refreshed) ---
$ws.Range("A3").Value = "Jairo Diaz"
$ws.Range("B3").Value = "diazja01"
$ws.Range("C3").Value = "September 12 2017"
$ws.Range("D3").Value = "Elbow"
$ws.Range("E3").Value = "Diaz is on the 60-day disabled list with inflammation in his left elbow and will miss the remainder of the season."

# --- Row 4: Carlos Gonzalez (new) ---
$ws.Range("A4").Value = "Carlos Gonzalez"
$ws.Range("B4").Value = "gonzaca01"
$ws.Range("C4").Value = "September 24 2017"
$ws.Range("D4").Value = "Shoulder"
$ws.Range("E4").Value = "Gonzalez sat out the last game due to a shoulder injury and his availability for Monday's tilt against the Marlins is undecided."

# Wrap text + taller rows for the three data rows (matches existing Injury.Details styling)
$ws.Range("E2:E4").WrapText = $true
$ws.Range("A2:E4").RowHeight = 30

# Two trailing blank rows, each with the wrap-text style parked on column E only
$ws.Range("E5").Value = ""
$ws.Range("E5").WrapText = $true
$ws.Range("E6").Value = ""
$ws.Range("E6").WrapText = $true

# Move the active selection the way the source workbook shows post-edit
$ws.Range("A14").Select() | Out-Null
